# Apply targeted username value edits on the "Users" sheet and update the
# active selection, matching the author's change (text-only formatting /
# value tweak, no structural changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Activate()

# Update usernames in column B (the " USER NAME" column)
$ws.Range("B5").Value = "Maria26"
$ws.Range("B2").Value = "BRUNO135"

# Move / set the active selection to B2 as in the final workbook state
$ws.Range("B2").Select()
